$wb = $excel.ActiveWorkbook

# --- Update "Users to Delete" sheet: new Database ID values (column D) ---
$wsUsers = $wb.Worksheets.Item("Users to Delete")

$wsUsers.Range("D2").Value = "3a87b84e-ad73-47ce-a66e-85738f511b12"
$wsUsers.Range("D3").Value = "48013d34-c952-42f2-a7d5-12d70ab56c71"
$wsUsers.Range("D4").Value = "33cbeaee-1d80-4176-8909-d9cc52bb113f"
$wsUsers.Range("D5").Value = "cc3e80ae-d970-4a7a-be85-8ce8bef43e6b"

# --- Update "Summary" sheet: Report Generated timestamp (B6) ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = "11/27/2025, 4:51:11 PM"
